$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row2 = @("ECs", "Timp1", "Cd63", "ECs", 2, 0.6666666666666666, 2.462258333333333, 7.386775, 0.04766847958679572, 0.04766847958679572, 3, 1, 88.76498299999999, 266.294949, 0.2826061427173424, 0.2826061427173424, 218.5623190999416, 1967.060871899475, 0.01347140514522471, 0.01347140514522471)
$row3 = @("ECs", "Timp1", "Cd63", "FAPs", 2, 0.6666666666666666, 2.462258333333333, 7.386775, 0.04766847958679572, 0.04766847958679572, 3, 1, 97.95392866666667, 293.861786, 0.3118615135035443, 0.3118615135035442, 241.1878771422389, 2170.69089428015, 0.01486596419035092, 0.01486596419035092)
$row4 = @("ECs", "Timp1", "Cd63", "sCs", 2, 0.6666666666666666, 2.462258333333333, 7.386775, 0.04766847958679572, 0.04766847958679572, 3, 1, 127.375404, 382.126212, 0.4055323437791133, 0.4055323437791133, 313.6311499607, 2822.6803496463, 0.01933111025122009, 0.01933111025122009)
$row5 = @("FAPs", "Timp1", "Cd63", "ECs", 3, 1, 27.97876166666667, 83.936285, 0.5416592610596597, 0.5416592610596597, 3, 1, 88.76498299999999, 266.294949, 0.2826061427173424, 0.2826061427173424, 2483.534303702718, 22351.80873332446, 0.1530762344351964, 0.1530762344351964)
$row6 = @("FAPs", "Timp1", "Cd63", "FAPs", 3, 1, 27.97876166666667, 83.936285, 0.5416592610596597, 0.5416592610596597, 3, 1, 97.95392866666667, 293.861786, 0.3118615135035443, 0.3118615135035442, 2740.629624478335, 24665.66662030501, 0.1689226769572769, 0.1689226769572769)
$row7 = @("FAPs", "Timp1", "Cd63", "sCs", 3, 1, 27.97876166666667, 83.936285, 0.5416592610596597, 0.5416592610596597, 3, 1, 127.375404, 382.126212, 0.4055323437791133, 0.4055323437791133, 3563.80607071138, 32074.25463640242, 0.2196603496671864, 0.2196603496671864)
$row8 = @("sCs", "Timp1", "Cd63", "ECs", 3, 1, 21.21278466666666, 63.63835399999999, 0.4106722593535446, 0.4106722593535446, 3, 1, 88.76498299999999, 266.294949, 0.2826061427173424, 0.2826061427173424, 1882.952470319327, 16946.57223287394, 0.1160585031369212, 0.1160585031369213)
$row9 = @("sCs", "Timp1", "Cd63", "FAPs", 3, 1, 21.21278466666666, 63.63835399999999, 0.4106722593535446, 0.4106722593535446, 3, 1, 97.95392866666667, 293.861786, 0.3118615135035443, 0.3118615135035442, 2077.875596060027, 18700.88036454024, 0.1280728723559165, 0.1280728723559165)
$row10 = @("sCs", "Timp1", "Cd63", "sCs", 3, 1, 21.21278466666666, 63.63835399999999, 0.4106722593535446, 0.4106722593535446, 3, 1, 127.375404, 382.126212, 0.4055323437791133, 0.4055323437791133, 2701.987016881672, 24317.88315193505, 0.1665408838607068, 0.1665408838607068)

$rows = @($row2, $row3, $row4, $row5, $row6, $row7, $row8, $row9, $row10)
$rowNums = @(2, 3, 4, 5, 6, 7, 8, 9, 10)

for ($ri = 0; $ri -lt $rows.Length; $ri++) {
  $rowVals = $rows[$ri]
  $rowNum = $rowNums[$ri]
  for ($ci = 0; $ci -lt $rowVals.Length; $ci++) {
    $ws.Cells.Item($rowNum, $ci + 1).Value = $rowVals[$ci]
  }
}
